# Update "想去人数" (F column) figures for the "展览" and "全部类型" sheets
# to reflect a newer scrape of the source data (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

# Map: sheet name -> @{ row = newValue }
$updates = @{
    "展览" = @{
        2  = 903
        5  = 1194
        6  = 74
        7  = 4388
        8  = 2618
        10 = 2538
        14 = 1665
        16 = 125
        17 = 117
        18 = 332
        26 = 561
        27 = 698
        28 = 108
        29 = 418
        30 = 51
        31 = 1623
        32 = 1071
        33 = 148
        34 = 22
        35 = 1165
        36 = 2078
        37 = 278
        40 = 87
        43 = 667
        44 = 1335
        45 = 112
        47 = 442
        48 = 72
    }
    "全部类型" = @{
        2  = 903
        3  = 1194
        5  = 74
        6  = 4388
        7  = 2618
        8  = 2538
        9  = 1665
        13 = 125
        14 = 117
        15 = 332
        22 = 561
        23 = 698
        24 = 108
        28 = 418
        29 = 1623
        30 = 1071
        31 = 148
        32 = 22
        34 = 2078
        35 = 278
        40 = 551
        41 = 87
        43 = 667
        44 = 1335
        46 = 112
        47 = 442
        48 = 72
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Range("F$row").Value = $rowsMap[$row]
    }
}
